# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# 1. Insert a new "Player Info" sheet before "ODI Batting" with ID/NAME/
#    BATTING_HAND/BOWL_STYLE columns for player 3902.
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE on the "ODI Batting" and
#    "ODI Bowling" sheets, replacing the full scorecard URL values with
#    just the numeric match code that used to be the URL's query param.

$wb = $excel.ActiveWorkbook

$battingBeforeInsert = $wb.Worksheets.Item("ODI Batting")

# --- 1. New "Player Info" sheet, inserted before "ODI Batting" ---------
$playerInfo = $wb.Worksheets.Add($battingBeforeInsert)
$playerInfo.Name = "Player Info"

# Re-resolve the other sheets by name now that the index positions have
# shifted because of the insert above.
$batting = $wb.Worksheets.Item("ODI Batting")
$bowling = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Cells.Item(1, 1).Value = "ID"
$playerInfo.Cells.Item(1, 2).Value = "NAME"
$playerInfo.Cells.Item(1, 3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1, 4).Value = "BOWL_STYLE"

# Copy the bold/bordered header formatting from the existing "ODI Batting"
# header row onto the new sheet's header row (values are untouched by a
# formats-only paste).
$batting.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$playerInfo.Cells.Item(2, 1).Value = "'3902"
$playerInfo.Cells.Item(2, 1).Style = "Normal"
$playerInfo.Cells.Item(2, 2).Value = "Nkruma Eljego Bonner"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Leg Break"

# --- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE ---------------------
$batting.Range("D1").Value = "MATCH_CODE"

$batting.Cells.Item(2, 4).Value = "'4443"
$batting.Cells.Item(3, 4).Value = "'4445"
$batting.Cells.Item(4, 4).Value = "'4447"
$batting.Cells.Item(5, 4).Value = "'4577"
$batting.Cells.Item(6, 4).Value = "'4580"
$batting.Cells.Item(7, 4).Value = "'4583"
$batting.Range("D2:D7").Style = "Normal"

# --- 2b. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE --------------------
$bowling.Range("B1").Value = "MATCH_CODE"

$bowling.Cells.Item(2, 2).Value = "'4443"
$bowling.Cells.Item(3, 2).Value = "'4580"
$bowling.Cells.Item(4, 2).Value = "'4583"
$bowling.Range("B2:B4").Style = "Normal"
